# New testcase of checkbox is added
#
# Adds a new "CheckBox" worksheet (after the existing "Home"/"InputForm"
# sheets), fills it with the checkbox test-case data, applies the
# bold/purple "success" style to the three "checked" label cells, and
# makes the new sheet the active tab - mirroring the authored commit.

$wb = $excel.ActiveWorkbook

# --- add the new sheet at the end of the workbook -------------------------
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "CheckBox"

# --- styled-cell font (bold, size 9, purple RGB 0x66,0x0E,0x7A) -----------
$successColor = 7999078  # OLE BGR encoding of RGB(0x66, 0x0E, 0x7A)

# Cells are written row-by-row, column A then column B, so that the shared-
# string table is populated in the same left-to-right / top-to-bottom order
# the original workbook uses.
$ws.Range("A1").Value = "checkbox"
$ws.Range("B1").Value = "checkBox"

$ws.Range("A2").Value = "checkbox1"
$b2 = $ws.Range("B2")
$b2.Value = "checkBox1"
$b2.Font.Bold = $true
$b2.Font.Size = 9
$b2.Font.Color = $successColor

$ws.Range("A3").Value = "checkbox2"
$b3 = $ws.Range("B3")
$b3.Value = "checkBox2"
$b3.Font.Bold = $true
$b3.Font.Size = 9
$b3.Font.Color = $successColor

$ws.Range("A4").Value = "checkbox3"
$ws.Range("B4").Value = "checkBox3"

$ws.Range("A5").Value = "checkbox4"
$b5 = $ws.Range("B5")
$b5.Value = "checkBox4"
$b5.Font.Bold = $true
$b5.Font.Size = 9
$b5.Font.Color = $successColor

$ws.Range("A6").Value = "status"
$ws.Range("B6").Value = "Check All"

$ws.Range("B7").Value = "Uncheck All"

$ws.Range("B9").Value = "Success - Check box is checked"

# --- selection / view state -------------------------------------------------
[void]$ws.Range("B10").Select()
